# edit.ps1 - apply the changes described by the diff to the active document
#
# 1. Insert two new paragraphs at the very top of the document body:
#       - "(new.docx)"  (Candara 10pt, en-GB, hanging-indent formatting that
#         matches the rest of the document)
#       - an empty paragraph with the same formatting
# 2. Move <w:lastRenderedPageBreak/> from just before the "Conducted " run to
#    just before the (first) "Automated " run.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Step 1: insert the two new paragraphs before the current first paragraph.
# ---------------------------------------------------------------------------

$newParasXml =
  '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:autoSpaceDE w:val="0"/>' +
      '<w:autoSpaceDN w:val="0"/>' +
      '<w:adjustRightInd w:val="0"/>' +
      '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
      '<w:ind w:left="2880" w:hanging="2880"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
        '<w:lang w:val="en-GB"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
        '<w:lang w:val="en-GB"/>' +
      '</w:rPr>' +
      '<w:t>(new.docx)</w:t>' +
    '</w:r>' +
  '</w:p>' +
  '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
      '<w:autoSpaceDE w:val="0"/>' +
      '<w:autoSpaceDN w:val="0"/>' +
      '<w:adjustRightInd w:val="0"/>' +
      '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
      '<w:ind w:left="2880" w:hanging="2880"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
        '<w:lang w:val="en-GB"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
  '</w:p>'

$topRange = $d.Paragraphs.First.Range
$topRange.Collapse(1)   # wdCollapseStart
$topRange.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# Step 2: move <w:lastRenderedPageBreak/> from "Conducted " to "Automated "
# ---------------------------------------------------------------------------
# The host's InsertXML inserts new paragraph-level content ahead of the
# target range rather than truly splicing in place, so a plain "insert at
# the same spot" ends up reordering the sibling runs. Instead: insert the
# replacement run *before* the target text, then delete the now-shifted
# original run text that followed it.

function Replace-RunXml($searchText, $replacementParaXml) {
    $rng = $d.Range(0, 0)
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Text = $searchText
    $find.Forward = $true
    $find.Wrap = 0
    $find.MatchCase = $true
    $find.Execute() | Out-Null

    $start = $rng.Start
    $end = $rng.End
    $len = $end - $start

    $insertion = $d.Range($start, $start)
    $insertion.InsertXML($replacementParaXml)

    $oldRange = $d.Range($start + $len, $end + $len)
    $oldRange.Delete()
}

# Remove the lastRenderedPageBreak sitting in front of "Conducted ".
$conductedXml =
  '<w:p ' + $wNs + '>' +
    '<w:r w:rsidRPr="00BF081F">' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
        '<w:lang w:val="en-GB"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">Conducted </w:t>' +
    '</w:r>' +
  '</w:p>'
Replace-RunXml "Conducted " $conductedXml

# Add the lastRenderedPageBreak in front of (the first) "Automated ".
$automatedXml =
  '<w:p ' + $wNs + '>' +
    '<w:r w:rsidRPr="00BF081F">' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Candara" w:hAnsi="Candara"/>' +
        '<w:sz w:val="20"/>' +
        '<w:szCs w:val="20"/>' +
        '<w:lang w:val="en-GB"/>' +
      '</w:rPr>' +
      '<w:lastRenderedPageBreak/>' +
      '<w:t xml:space="preserve">Automated </w:t>' +
    '</w:r>' +
  '</w:p>'
Replace-RunXml "Automated " $automatedXml
